# Vancouver Home Pass Types — add a "Match ID" column in front of the
# existing data (everything shifts one column to the right, A:V -> B:W)
# and populate the new column A with the match id (13) for every data
# row, plus the "Match ID" header in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts all existing
# cells/styles/merged-cell refs one column to the right automatically.
$ws.Columns("A").Insert()

# New header label (row 2 holds the real column headers).
$ws.Range("A2").Value = "Match ID"

# Populate the match id for all player rows (4-19) and the hidden
# totals row (20). Row 3 (hidden spacer) stays blank.
$ws.Range("A4:A19").Value = 13
$ws.Range("A20").Value = 13

# Match the bold header styling used elsewhere in row 2, but without a
# border (A20's cell keeps the default/no style).
$ws.Range("A2:A19").Font.Bold = $true

# Re-fit the previously-blank hidden rows so they don't pick up a
# stray explicit row height from the new writes.
$ws.Rows("1").AutoFit()
$ws.Rows("3").AutoFit()
$ws.Rows("20").AutoFit()

# Restore the on-screen selection to the newly added column's data
# range, as in the saved workbook.
$ws.Range("A2:A19").Select()
